# Scheduled-runner update: refresh market-derived profit figures across the
# Sheets workbook (currentAveragePrice / LevePrice / LeveProfit columns,
# i.e. columns H..N). Identifying columns A..G are left untouched.

$wb = $excel.ActiveWorkbook

function SetCell($SheetName, $Ref, $Val) {
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($Ref).Value = $Val
}

function ClearCell($SheetName, $Ref) {
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($Ref).Value = ""
}

# ---------------------------------------------------------------- ALC ----
SetCell "ALC" "H43" 1940.6666
SetCell "ALC" "J43" 1218.8
SetCell "ALC" "L43" 1218.8
SetCell "ALC" "N43" -1356.8

SetCell "ALC" "H98" 2885
SetCell "ALC" "I98" 2885
SetCell "ALC" "K98" 2885
SetCell "ALC" "M98" -1387

SetCell "ALC" "H106" 1824.75
SetCell "ALC" "I106" 1824.75
SetCell "ALC" "K106" 1824.75
SetCell "ALC" "M106" -1193.75

SetCell "ALC" "H122" 2885
SetCell "ALC" "I122" 2885
SetCell "ALC" "K122" 8655
SetCell "ALC" "M122" -6205

SetCell "ALC" "H127" 1951.4814
SetCell "ALC" "I127" 681.1111
SetCell "ALC" "J127" 2586.6667
SetCell "ALC" "K127" 2043.3333
SetCell "ALC" "L127" 7760.000100000001
SetCell "ALC" "M127" 2916.6667
SetCell "ALC" "N127" -17680.0001

SetCell "ALC" "H132" 1494203.8
SetCell "ALC" "I132" 1708.8036
SetCell "ALC" "J132" 9092360
SetCell "ALC" "K132" 5126.4108
SetCell "ALC" "L132" 27277080
SetCell "ALC" "M132" -2596.4108
SetCell "ALC" "N132" -27282140

# ---------------------------------------------------------------- ARM ----
SetCell "ARM" "H32" 6008.57
SetCell "ARM" "I32" 6008.57
SetCell "ARM" "J32" 0
SetCell "ARM" "K32" 6008.57
SetCell "ARM" "L32" 0
SetCell "ARM" "M32" -5721.57
ClearCell "ARM" "N32"

SetCell "ARM" "H61" 1674.2667
SetCell "ARM" "I61" 1658.1428
SetCell "ARM" "K61" 1658.1428
SetCell "ARM" "M61" -1446.1428

SetCell "ARM" "H63" 2061.8125
SetCell "ARM" "I63" 2069.2666
SetCell "ARM" "J63" 1950
SetCell "ARM" "K63" 2069.2666
SetCell "ARM" "L63" 1950
SetCell "ARM" "M63" -1383.2666
SetCell "ARM" "N63" -3322

SetCell "ARM" "H66" 2061.8125
SetCell "ARM" "I66" 2069.2666
SetCell "ARM" "J66" 1950
SetCell "ARM" "K66" 10346.333
SetCell "ARM" "L66" 9750
SetCell "ARM" "M66" -6914.332999999999
SetCell "ARM" "N66" -16614

SetCell "ARM" "H132" 1101.0303
SetCell "ARM" "I132" 1101.0303
SetCell "ARM" "J132" 0
SetCell "ARM" "K132" 3303.0909
SetCell "ARM" "L132" 0
SetCell "ARM" "M132" -773.0908999999997
ClearCell "ARM" "N132"

SetCell "ARM" "H136" 1674.2667
SetCell "ARM" "I136" 1658.1428
SetCell "ARM" "K136" 4974.428400000001
SetCell "ARM" "M136" -2424.428400000001

# ---------------------------------------------------------------- BSM ----
SetCell "BSM" "H134" 1688
SetCell "BSM" "I134" 1474.8718
SetCell "BSM" "J134" 10000
SetCell "BSM" "K134" 4424.6154
SetCell "BSM" "L134" 30000
SetCell "BSM" "M134" -1889.6154
SetCell "BSM" "N134" -35070

# ---------------------------------------------------------------- CUL ----
SetCell "CUL" "H113" 697.8148
SetCell "CUL" "J113" 620.7222
SetCell "CUL" "L113" 1862.1666
SetCell "CUL" "N113" -6202.1666

SetCell "CUL" "H133" 3500
SetCell "CUL" "I133" 3000
SetCell "CUL" "J133" 4000
SetCell "CUL" "K133" 9000
SetCell "CUL" "L133" 12000
SetCell "CUL" "M133" -3940
SetCell "CUL" "N133" -22120

SetCell "CUL" "H134" 4408.7085
SetCell "CUL" "I134" 2907.2666
SetCell "CUL" "J134" 6911.1113
SetCell "CUL" "K134" 8721.799800000001
SetCell "CUL" "L134" 20733.3339
SetCell "CUL" "M134" -3651.799800000001
SetCell "CUL" "N134" -30873.3339

SetCell "CUL" "H138" 2321
SetCell "CUL" "I138" 1370
SetCell "CUL" "J138" 5491
SetCell "CUL" "K138" 4110
SetCell "CUL" "L138" 16473
SetCell "CUL" "M138" 1030
SetCell "CUL" "N138" -26753

SetCell "CUL" "H139" 33708.066
SetCell "CUL" "I139" 42664.582
SetCell "CUL" "J139" 3000
SetCell "CUL" "K139" 127993.746
SetCell "CUL" "L139" 9000
SetCell "CUL" "M139" -122853.746
SetCell "CUL" "N139" -19280

# ---------------------------------------------------------------- LTW ----
SetCell "LTW" "H7" 2696.8333
SetCell "LTW" "I7" 2760.182
SetCell "LTW" "J7" 2000
SetCell "LTW" "K7" 2760.182
SetCell "LTW" "L7" 2000
SetCell "LTW" "M7" -2648.182
SetCell "LTW" "N7" -2224

SetCell "LTW" "H40" 2138.8333
SetCell "LTW" "I40" 2031.1875
SetCell "LTW" "K40" 2031.1875
SetCell "LTW" "M40" -1895.1875

SetCell "LTW" "H61" 1430.3948
SetCell "LTW" "I61" 1381.9642
SetCell "LTW" "K61" 1381.9642
SetCell "LTW" "M61" -1179.9642

SetCell "LTW" "H109" 21800
SetCell "LTW" "J109" 21800
SetCell "LTW" "L109" 21800
SetCell "LTW" "N109" -24574

SetCell "LTW" "H113" 1430.3948
SetCell "LTW" "I113" 1381.9642
SetCell "LTW" "K113" 1381.9642
SetCell "LTW" "M113" 788.0358000000001

SetCell "LTW" "H126" 2696.8333
SetCell "LTW" "I126" 2760.182
SetCell "LTW" "J126" 2000
SetCell "LTW" "K126" 8280.545999999998
SetCell "LTW" "L126" 6000
SetCell "LTW" "M126" -5810.545999999998
SetCell "LTW" "N126" -10940

SetCell "LTW" "H132" 2473.7083
SetCell "LTW" "I132" 1578.28
SetCell "LTW" "J132" 3447
SetCell "LTW" "K132" 4734.84
SetCell "LTW" "L132" 10341
SetCell "LTW" "M132" -2204.84
SetCell "LTW" "N132" -15401

# ---------------------------------------------------------------- WVR ----
SetCell "WVR" "H43" 35500
SetCell "WVR" "I43" 35000
SetCell "WVR" "J43" 36000
SetCell "WVR" "K43" 35000
SetCell "WVR" "L43" 36000
SetCell "WVR" "M43" -34851
SetCell "WVR" "N43" -36298

SetCell "WVR" "H100" 1294
SetCell "WVR" "I100" 1150.8572
SetCell "WVR" "K100" 2301.7144
SetCell "WVR" "M100" -1760.7144

SetCell "WVR" "H126" 1516.6666
SetCell "WVR" "I126" 1516.6666
SetCell "WVR" "J126" 0
SetCell "WVR" "K126" 4549.9998
SetCell "WVR" "L126" 0
SetCell "WVR" "M126" -2079.9998
ClearCell "WVR" "N126"

SetCell "WVR" "H132" 832.57776
SetCell "WVR" "I132" 832.57776
SetCell "WVR" "J132" 0
SetCell "WVR" "K132" 2497.73328
SetCell "WVR" "L132" 0
SetCell "WVR" "M132" 32.26672000000008
ClearCell "WVR" "N132"
